# Auto-generated edit script applying the Phantom_Profits diff
# Updates numeric cells (H:N) across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6533.222
$ws.Range("I51").Value = 6239.8
$ws.Range("J51").Value = 6900
$ws.Range("K51").Value = 6239.8
$ws.Range("L51").Value = 6900
$ws.Range("M51").Value = -5755.8
$ws.Range("N51").Value = -7868
$ws.Range("H64").Value = 5466
$ws.Range("I64").Value = 4599
$ws.Range("K64").Value = 4599
$ws.Range("M64").Value = -4351
$ws.Range("H67").Value = 5466
$ws.Range("I67").Value = 4599
$ws.Range("K67").Value = 4599
$ws.Range("M67").Value = -3741
$ws.Range("H76").Value = 2994.2
$ws.Range("I76").Value = 2942.75
$ws.Range("K76").Value = 2942.75
$ws.Range("M76").Value = -2627.75
$ws.Range("H79").Value = 2994.2
$ws.Range("I79").Value = 2942.75
$ws.Range("K79").Value = 2942.75
$ws.Range("M79").Value = -1850.75
$ws.Range("H95").Value = 31049.5
$ws.Range("J95").Value = 31049.5
$ws.Range("L95").Value = 31049.5
$ws.Range("N95").Value = -36541.5
$ws.Range("H98").Value = 450.7143
$ws.Range("I98").Value = 450.7143
$ws.Range("K98").Value = 450.7143
$ws.Range("M98").Value = 1047.2857
$ws.Range("H122").Value = 450.7143
$ws.Range("I122").Value = 450.7143
$ws.Range("K122").Value = 1352.1429
$ws.Range("M122").Value = 1097.8571
$ws.Range("H125").Value = 2745
$ws.Range("I125").Value = 3017
$ws.Range("J125").Value = 1929
$ws.Range("K125").Value = 27153
$ws.Range("L125").Value = 17361
$ws.Range("M125").Value = -24693
$ws.Range("N125").Value = -22281
$ws.Range("H141").Value = 18023
$ws.Range("I141").Value = 17998.5
$ws.Range("K141").Value = 53995.5
$ws.Range("M141").Value = -48815.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1126.6111
$ws.Range("I32").Value = 1161.1471
$ws.Range("K32").Value = 1161.1471
$ws.Range("M32").Value = -874.1470999999999
$ws.Range("H122").Value = 907.55554
$ws.Range("I122").Value = 907.55554
$ws.Range("K122").Value = 2722.66662
$ws.Range("M122").Value = -272.66662
$ws.Range("H131").Value = 20000
$ws.Range("J131").Value = 20000
$ws.Range("L131").Value = 20000
$ws.Range("N131").Value = -30080

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4547.5625
$ws.Range("I86").Value = 5050.5
$ws.Range("J86").Value = 3709.3333
$ws.Range("K86").Value = 5050.5
$ws.Range("L86").Value = 3709.3333
$ws.Range("M86").Value = -3927.5
$ws.Range("N86").Value = -5955.3333
$ws.Range("H89").Value = 4547.5625
$ws.Range("I89").Value = 5050.5
$ws.Range("J89").Value = 3709.3333
$ws.Range("K89").Value = 25252.5
$ws.Range("L89").Value = 18546.6665
$ws.Range("M89").Value = -19636.5
$ws.Range("N89").Value = -29778.6665
$ws.Range("H105").Value = 2891.2727
$ws.Range("I105").Value = 1861.1818
$ws.Range("K105").Value = 1861.1818
$ws.Range("M105").Value = -114.1818000000001
$ws.Range("H134").Value = 3529.818
$ws.Range("I134").Value = 3870.111
$ws.Range("K134").Value = 11610.333
$ws.Range("M134").Value = -9075.332999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 784.875
$ws.Range("I22").Value = 723.375
$ws.Range("K22").Value = 723.375
$ws.Range("M22").Value = -373.375
$ws.Range("H107").Value = 1943.0526
$ws.Range("I107").Value = 1846.3572
$ws.Range("J107").Value = 2213.8
$ws.Range("K107").Value = 1846.3572
$ws.Range("L107").Value = 2213.8
$ws.Range("M107").Value = 73.64280000000008
$ws.Range("N107").Value = -6053.8
$ws.Range("H134").Value = 3602.8
$ws.Range("I134").Value = 3602.8
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 10808.4
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -8273.400000000001
$ws.Range("N134").ClearContents()
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 38.25
$ws.Range("I6").Value = 38.25
$ws.Range("K6").Value = 114.75
$ws.Range("M6").Value = -1.75
$ws.Range("H7").Value = 326.8
$ws.Range("I7").Value = 415.66666
$ws.Range("K7").Value = 1246.99998
$ws.Range("M7").Value = -1134.99998
$ws.Range("H131").Value = 1865
$ws.Range("I131").Value = 1230
$ws.Range("J131").Value = 2500
$ws.Range("K131").Value = 3690
$ws.Range("L131").Value = 7500
$ws.Range("M131").Value = 1350
$ws.Range("N131").Value = -17580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 42018
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H80").Value = 873.4286
$ws.Range("I80").Value = 742.8
$ws.Range("J80").Value = 1200
$ws.Range("K80").Value = 742.8
$ws.Range("L80").Value = 1200
$ws.Range("M80").Value = 255.2
$ws.Range("N80").Value = -3196
$ws.Range("H83").Value = 873.4286
$ws.Range("I83").Value = 742.8
$ws.Range("J83").Value = 1200
$ws.Range("K83").Value = 3714
$ws.Range("L83").Value = 6000
$ws.Range("M83").Value = 1278
$ws.Range("N83").Value = -15984
$ws.Range("H102").Value = 2238.3333
$ws.Range("I102").Value = 2286
$ws.Range("K102").Value = 2286
$ws.Range("M102").Value = -664
$ws.Range("H113").Value = 2344
$ws.Range("I113").Value = 1702
$ws.Range("J113").Value = 3200
$ws.Range("K113").Value = 1702
$ws.Range("L113").Value = 3200
$ws.Range("M113").Value = 468
$ws.Range("N113").Value = -7540
$ws.Range("H128").Value = 49250
$ws.Range("J128").Value = 30000
$ws.Range("L128").Value = 30000
$ws.Range("N128").Value = -39960
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4999
$ws.Range("I7").Value = 4999
$ws.Range("K7").Value = 4999
$ws.Range("M7").Value = -4887
$ws.Range("H22").Value = 1047.75
$ws.Range("I22").Value = 895.5
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 895.5
$ws.Range("L22").Value = 1200
$ws.Range("M22").Value = -600.5
$ws.Range("N22").Value = -1790
$ws.Range("H27").Value = 1047.75
$ws.Range("I27").Value = 895.5
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 895.5
$ws.Range("L27").Value = 1200
$ws.Range("M27").Value = -788.5
$ws.Range("N27").Value = -1414
$ws.Range("H55").Value = 1119.2307
$ws.Range("I55").Value = 1068.5
$ws.Range("J55").Value = 1200.4
$ws.Range("K55").Value = 1068.5
$ws.Range("L55").Value = 1200.4
$ws.Range("M55").Value = -895.5
$ws.Range("N55").Value = -1546.4
$ws.Range("H82").Value = 5255.778
$ws.Range("I82").Value = 3866.6667
$ws.Range("K82").Value = 3866.6667
$ws.Range("M82").Value = -3505.6667
$ws.Range("H85").Value = 5255.778
$ws.Range("I85").Value = 3866.6667
$ws.Range("K85").Value = 3866.6667
$ws.Range("M85").Value = -2618.6667
$ws.Range("H93").Value = 1427
$ws.Range("J93").Value = 1448.3334
$ws.Range("L93").Value = 1448.3334
$ws.Range("N93").Value = -3944.3334
$ws.Range("H122").Value = 20999.6
$ws.Range("J122").Value = 19999.5
$ws.Range("L122").Value = 59998.5
$ws.Range("N122").Value = -64898.5
$ws.Range("H126").Value = 4999
$ws.Range("I126").Value = 4999
$ws.Range("K126").Value = 14997
$ws.Range("M126").Value = -12527
$ws.Range("H132").Value = 2876.25
$ws.Range("I132").Value = 2833.8333
$ws.Range("K132").Value = 8501.499899999999
$ws.Range("M132").Value = -5971.499899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2482.8333
$ws.Range("I126").Value = 2579.4
$ws.Range("K126").Value = 7738.200000000001
$ws.Range("M126").Value = -5268.200000000001
$ws.Range("H130").Value = 29999
$ws.Range("J130").Value = 29999
$ws.Range("L130").Value = 29999
$ws.Range("N130").Value = -40039
$ws.Range("H136").Value = 5829
$ws.Range("I136").Value = 6200.6
$ws.Range("K136").Value = 18601.8
$ws.Range("M136").Value = -16051.8
